$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellsToClear = @("I2","J2","K2","G3","H3","I3","J3","K3","I4","J4","K4","G5","H5","I5","J5","K5","I6","J6","K6","G7","H7","I7","J7","K7","I8","J8","K8","G9","H9","I9","J9","K9","I10","J10","K10","G11","H11","I11","J11","K11","I12","J12","K12","K13","I14","J14","K14","K15","I16","J16","K16","K17","J18","K18","I19","J19","K19","K21","J22","K22","I23","J23","K23","K25","J26","K26","I27","J27","K27","K29","J30","K30","I31","J31","K31","K33","J34","K34","I35","J35","K35","K37","J38","K38","I39","J39")

foreach ($addr in $cellsToClear) {
    $ws.Range($addr).ClearContents()
}
